$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.1366576124457447;  C = 0.775372048323969;   D = 1.844561190472676;  E = 1.358146233095934;  F = 1.364436723157444;  G = 52 }
    3  = @{ B = 0.100644678693309;   C = 0.5870089618847463;  D = 0.9743935215277069; E = 0.9871137328229745; F = 0.9917406180962824; G = 51 }
    4  = @{ B = 0.09787366134911048; C = 0.6307949273637553;  D = 1.168272060271323;  E = 1.080866347089835;  F = 1.087354403430002;  G = 50 }
    5  = @{ B = 0.08154868731045999; C = 0.713998770293546;   D = 1.430771530282745;  E = 1.196148623826799;  F = 1.208190537716507;  G = 41 }
    6  = @{ B = 0.1341520404658785;  C = 0.7691896598824904;  D = 1.579558871789514;  E = 1.256805025367703;  F = 1.270281138917405;  G = 31 }
    7  = @{ B = 0.1268598637344239;  C = 0.7813345819487484;  D = 1.592716119073123;  E = 1.262028573001865;  F = 1.277861748157917;  G = 29 }
    8  = @{ B = 0.09675479195317427; C = 0.8108202300352334;  D = 1.726960037846746;  E = 1.314138515471922;  F = 1.335537356534567;  G = 27 }
    9  = @{ B = 0.1406368155206148;  C = 0.8950670472051415;  D = 2.222009981461377;  E = 1.490640795584697;  F = 1.524656502007553;  G = 19 }
    10 = @{ B = -0.3312026136842777; C = 0.8013776111497077;  D = 1.293401115164856;  E = 1.137277941035021;  F = 1.136360671227658;  G = 12 }
    11 = @{ B = 0.4641315885571715;  C = 0.7592743174351484;  D = 1.498317365315948;  E = 1.224057745907418;  F = 1.266342782296397;  G = 5 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
